$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the contribution notes for each team member (percentages revised).
$ws.Range("I11").Value = "UI/Functions 22%"
$ws.Range("I12").Value = "UI/Functions 22%"
$ws.Range("I13").Value = "DB/Functions/ Manager 22%"
$ws.Range("I14").Value = "DB/Manager/UI/Functions 30%"

# Move the active selection to I22 (matches the saved view state).
$ws.Range("I22").Select()
